# Add data for 2022-12-26 (shifts the "through Dec 17" cutoff to "through Dec 18")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and update the header label / shared string text
$ws.Name = "Through 2022-12-18"
$ws.Range("B1").Value = "December 2022 (through December 18)"

# Update existing counts that increased by one additional day of data
$ws.Range("BV4").Value = 4
$ws.Range("B5").Value = 3
$ws.Range("BJ10").Value = 3
$ws.Range("B12").Value = 3
$ws.Range("Z14").Value = 6
$ws.Range("AL14").Value = 3
$ws.Range("AX14").Value = 9
$ws.Range("N15").Value = 6
$ws.Range("AL20").Value = 4
$ws.Range("AX24").Value = 3
$ws.Range("N35").Value = 4
$ws.Range("Z41").Value = 4
$ws.Range("N55").Value = 3

# Newly populated cells (previously blank)
$ws.Range("AX5").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("B8").Value = 2
$ws.Range("BJ15").Value = 1
$ws.Range("CH33").Value = 1
$ws.Range("BV35").Value = 1
$ws.Range("Z38").Value = 1
$ws.Range("CH52").Value = 1
$ws.Range("BJ57").Value = 1
$ws.Range("BV59").Value = 1
